$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextCell $ws "A2" "PESCO_30262230000900.pdf"
Set-TextCell $ws "B2" "PESCO"
Set-TextCell $ws "C2" "30262230000900"
Set-TextCell $ws "D2" "30262230000900"
Set-TextCell $ws "E2" "489"
Set-TextCell $ws "F2" "489"
Set-TextCell $ws "G2" "A-2c(06)T"
Set-TextCell $ws "H2" "SEP-2025"
Set-TextCell $ws "I2" "966456"
Set-TextCell $ws "J2" "568"
Set-TextCell $ws "K2" "13173"
Set-TextCell $ws "L2" "176"
Set-TextCell $ws "M2" "3062"
Set-TextCell $ws "N2" "157"
Set-TextCell $ws "O2" "157"
Set-TextCell $ws "Q2" "35.1500"
Set-TextCell $ws "R2" "43.8200"
Set-TextCell $ws "S2" "D:\BILLs App\BillWebApp\BillWebApp\downloads\PESCO_30262230000900.pdf"

# Row 3
Set-TextCell $ws "A3" "PESCO_30265520001100.pdf"
Set-TextCell $ws "B3" "PESCO"
Set-TextCell $ws "C3" "30265520001100"
Set-TextCell $ws "H3" "NOT FOUND"
Set-TextCell $ws "S3" "D:\BILLs App\BillWebApp\BillWebApp\downloads\PESCO_30265520001100.pdf"

# Row 4
Set-TextCell $ws "A4" "PESCO_30266220052100.pdf"
Set-TextCell $ws "B4" "PESCO"
Set-TextCell $ws "C4" "30266220052100"
Set-TextCell $ws "D4" "30266220052100"
Set-TextCell $ws "E4" "212"
Set-TextCell $ws "F4" "68"
Set-TextCell $ws "G4" "B2b(12)T"
Set-TextCell $ws "H4" "SEP-2025"
Set-TextCell $ws "I4" "96073"
Set-TextCell $ws "J4" "120"
Set-TextCell $ws "K4" "867"
Set-TextCell $ws "L4" "18"
Set-TextCell $ws "M4" "349"
Set-TextCell $ws "N4" "27"
Set-TextCell $ws "O4" "27"
Set-TextCell $ws "Q4" "27.4100"
Set-TextCell $ws "R4" "36.6800"
Set-TextCell $ws "S4" "D:\BILLs App\BillWebApp\BillWebApp\downloads\PESCO_30266220052100.pdf"

# Row 5
Set-TextCell $ws "A5" "PESCO_30262430001100.pdf"
Set-TextCell $ws "B5" "PESCO"
Set-TextCell $ws "C5" "30262430001100"
Set-TextCell $ws "D5" "30262430001100"
Set-TextCell $ws "E5" "212"
Set-TextCell $ws "F5" "68"
Set-TextCell $ws "G5" "B2b(12)T"
Set-TextCell $ws "H5" "SEP-2025"
Set-TextCell $ws "I5" "158694"
Set-TextCell $ws "J5" "360"
Set-TextCell $ws "K5" "2360"
Set-TextCell $ws "L5" "160"
Set-TextCell $ws "M5" "1080"
Set-TextCell $ws "N5" "40"
Set-TextCell $ws "O5" "40"
Set-TextCell $ws "Q5" "27.4100"
Set-TextCell $ws "R5" "36.6800"
Set-TextCell $ws "S5" "D:\BILLs App\BillWebApp\BillWebApp\downloads\PESCO_30262430001100.pdf"

# Row 6
Set-TextCell $ws "A6" "PESCO_30262340005030.pdf"
Set-TextCell $ws "B6" "PESCO"
Set-TextCell $ws "C6" "30262340005030"
Set-TextCell $ws "D6" "30262340005030"
Set-TextCell $ws "E6" "212"
Set-TextCell $ws "F6" "36"
Set-TextCell $ws "G6" "B2b(12)T"
Set-TextCell $ws "H6" "SEP-2025"
Set-TextCell $ws "I6" "76997"
Set-TextCell $ws "J6" "98"
Set-TextCell $ws "K6" "979"
Set-TextCell $ws "L6" "66"
Set-TextCell $ws "M6" "386"
Set-TextCell $ws "N6" "25"
Set-TextCell $ws "O6" "25"
Set-TextCell $ws "Q6" "27.4100"
Set-TextCell $ws "R6" "36.6800"
Set-TextCell $ws "S6" "D:\BILLs App\BillWebApp\BillWebApp\downloads\PESCO_30262340005030.pdf"
